$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 11 ("E2E IOAM Indicator Label Allocation Methods"): add 6pt
# space-before to every bulleted paragraph in the content placeholder.
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$tr11 = $s11.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
for ($i = 1; $i -le $tr11.Paragraphs().Count; $i++) {
    $tr11.Paragraphs($i).ParagraphFormat.SpaceBefore = 6
}

# ---------------------------------------------------------------------------
# Slide 13 ("E2E IOAM Procedure"): tighten line spacing from 21.4pt to
# 20.4pt on every bullet in the content placeholder.
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$tr13 = $s13.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
for ($i = 1; $i -le $tr13.Paragraphs().Count; $i++) {
    $tr13.Paragraphs($i).ParagraphFormat.SpaceWithin = 20.4
}

# ---------------------------------------------------------------------------
# Slide 14: title "Example 1 - SR-MPLS Encapsulation..." -> "Example - ..."
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$s14.Shapes.Item("Title 1").TextFrame.TextRange.Runs(1).Text = "Example - SR-MPLS Encapsulation with IOAM Data Fields"

# ---------------------------------------------------------------------------
# Slide 16 ("HbH IOAM Indicator Label Allocation Methods"): add 6pt
# space-before to every bullet in the content placeholder.
# ---------------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$tr16 = $s16.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
for ($i = 1; $i -le $tr16.Paragraphs().Count; $i++) {
    $tr16.Paragraphs($i).ParagraphFormat.SpaceBefore = 6
}

# ---------------------------------------------------------------------------
# Slide 20: tighten line spacing from 22.2pt to 21.2pt.
# ---------------------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$tr20 = $s20.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
for ($i = 1; $i -le $tr20.Paragraphs().Count; $i++) {
    $tr20.Paragraphs($i).ParagraphFormat.SpaceWithin = 21.2
}

# ---------------------------------------------------------------------------
# Slide 7 ("IOAM G-ACh Header"): two small text tweaks.
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
$tr7.Paragraphs(3).Characters(1, 11).Text = "Note: GAL with G-"
$tr7.Paragraphs(4).Runs(1).Text = "Block Number is used to: "

# ---------------------------------------------------------------------------
# Slide 27: reposition "Rectangle 6" callout box.
# ---------------------------------------------------------------------------
$s27 = $p.Slides.Item(27)
$rect6 = $s27.Shapes.Item("Rectangle 6")
$rect6.Left = -1.5
$rect6.Top = 320.76515905511815

# ---------------------------------------------------------------------------
# Slide 8 ("IOAM Indicator Labels"): tighten line spacing from 23.2pt to
# 22.2pt.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item("Content Placeholder 2").TextFrame.TextRange
for ($i = 1; $i -le $tr8.Paragraphs().Count; $i++) {
    $tr8.Paragraphs($i).ParagraphFormat.SpaceWithin = 22.2
}
